$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for which data values are present in each row (skips A=index, G constant=1, H empty)
$cols = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

$data = @(
    @(2, 1.02, 1.089903957704872, 1.092945431346923, 1.092264471883205, 1.104236304247218, 1.073478190317613, 1.094737009875696, 1.095586576172583, 1.094907354187842, 1.106848976647145, 1.096291661997809),
    @(3, 1.02, 1.091259431272322, 1.094071114200453, 1.093478896407954, 1.105485783000437, 1.074014326944082, 1.095755481156932, 1.096532108786282, 1.095941292596534, 1.10792007710412, 1.097311579625085),
    @(4, 1.02, 1.092135457746356, 1.09479854640817, 1.094263457120968, 1.106293408140916, 1.074359446319635, 1.096412903088916, 1.097142377704407, 1.096608490950811, 1.10861170678259, 1.09796993517162),
    @(5, 1.02, 1.092503491361569, 1.095104132493929, 1.094592990505612, 1.106632729038955, 1.074504106454917, 1.096688904537662, 1.097398565999893, 1.096888547108778, 1.108902125833917, 1.098246328574081),
    @(6, 1.02, 1.092565271426049, 1.095155428563302, 1.094648303374915, 1.106689690617997, 1.074528370488687, 1.096735224291183, 1.097441559612101, 1.096935544429046, 1.1089508685165, 1.098292714106965),
    @(7, 1.02, 1.092140376397704, 1.094802630549852, 1.094267861519947, 1.106297942962945, 1.074361380954749, 1.096416592517441, 1.097145802349874, 1.096612234775584, 1.10861558871575, 1.097973629839557),
    @(8, 1.02, 1.090362267792081, 1.09332606192476, 1.092675154266931, 1.104658754142405, 1.073659753536407, 1.095081539882038, 1.095906447054545, 1.095257159510347, 1.107211261219926, 1.096636681276279),
    @(9, 1.02, 1.0872207139044, 1.090716655826773, 1.089858820041049, 1.101763449015367, 1.072409538977735, 1.092716612321294, 1.093710497431337, 1.092855164737272, 1.104725435361152, 1.094268395247137),
    @(10, 1.02, 1.085120462121774, 1.088971791965699, 1.087974430314145, 1.099828399401921, 1.071566610870064, 1.091131431768333, 1.09223822488046, 1.091244060018079, 1.103060444250978, 1.092680963555968),
    @(11, 1.02, 1.084209566758756, 1.08821495060247, 1.087156789956307, 1.098989297123908, 1.071199341472583, 1.090442948720099, 1.091598698027644, 1.090544061076366, 1.102337589232726, 1.091991502782836),
    @(12, 1.02, 1.083870992503704, 1.087933625872079, 1.086852823357983, 1.09867743010501, 1.071062576641601, 1.09018689702362, 1.09136084123704, 1.09028368796714, 1.102068798269429, 1.091735087463579),
    @(13, 1.02, 1.083943628205843, 1.087993980107758, 1.086918037007696, 1.098744335182454, 1.071091928800362, 1.090241835464654, 1.091411876331572, 1.090339555356236, 1.102126468058955, 1.091790103923501),
    @(14, 1.02, 1.08418158477365, 1.088191700320277, 1.08713166928661, 1.098963521943076, 1.071188043494807, 1.090421789929303, 1.091579043018856, 1.090522545973144, 1.102315376838761, 1.09197031394413),
    @(15, 1.02, 1.084328167389661, 1.08831349561886, 1.08726326076955, 1.099098545206616, 1.071247217183595, 1.09053262348226, 1.091681998990668, 1.090635244375616, 1.102431731152525, 1.092081304893446),
    @(16, 1.02, 1.085180883717644, 1.089021993201109, 1.088028658431676, 1.099884061888063, 1.07159093714951, 1.091177079713476, 1.092280625197923, 1.091290466092582, 1.103108377279617, 1.09272667632643),
    @(17, 1.02, 1.08571537173567, 1.089466062946872, 1.088508316878153, 1.100376467587558, 1.071805932370926, 1.091580767498325, 1.092655583073697, 1.09170082881915, 1.103532306867253, 1.093130937394236),
    @(18, 1.02, 1.086026987884094, 1.089724955809614, 1.088787931093103, 1.100663562790338, 1.071931116031274, 1.091816030463681, 1.092874094735166, 1.091939957146878, 1.103779394531731, 1.093366534459988),
    @(19, 1.02, 1.086133217086921, 1.089813210425656, 1.088883244914659, 1.100761435184353, 1.071973763288686, 1.091896215049177, 1.092948568625542, 1.092021454976064, 1.103863614087118, 1.093446832916793),
    @(20, 1.02, 1.085658040927755, 1.089418431479701, 1.0884568709066, 1.10032364922436, 1.071782888133676, 1.091537476473825, 1.092615373834828, 1.091656824536439, 1.103486842227362, 1.093087584891516),
    @(21, 1.02, 1.084111518813895, 1.088133482205269, 1.087068767083733, 1.098898982108179, 1.071159749647848, 1.090368806651641, 1.091529825132431, 1.090468669838053, 1.102259755965249, 1.091917255424137),
    @(22, 1.02, 1.083137841741307, 1.087324423459032, 1.086194512606173, 1.098002152595803, 1.070765962519898, 1.089632172907317, 1.090845513295484, 1.089719531082494, 1.101486555557192, 1.091179575575403),
    @(23, 1.02, 1.083654133228514, 1.08775343220755, 1.086658115192265, 1.098477683329588, 1.070974906552145, 1.090022852761285, 1.09120845044881, 1.090116864246335, 1.101896604732625, 1.091570810239575),
    @(24, 1.02, 1.085683946689453, 1.089439954477024, 1.088480117628315, 1.100347515929516, 1.071793301510359, 1.091557038446044, 1.092633543261138, 1.091676708885897, 1.103507386307216, 1.09310717464398),
    @(25, 1.02, 1.088033895453, 1.091392161095333, 1.090588095437072, 1.102512790994027, 1.072734405415389, 1.093329495063093, 1.094279651867689, 1.093477843056722, 1.105369434246184, 1.094882148352729),
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$colLetter$row").Value = $value
    }
}
